# Fruta / hortaliza, semanal
# Insert this week's two new price records (date 44511) at the top of the
# data block (row 80), pushing the existing rows (old 80-103) down to
# (82-105). The sheet keeps growing week over week with the newest
# observations inserted right after the header/earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows; Insert() shifts rows 80:103 down to 82:105
# and carries the row-80 formatting (incl. the date style on column D) along.
$ws.Rows("80:81").Insert()

# New record 1 -> row 80
$ws.Cells.Item(80, 1).Value  = 4
$ws.Cells.Item(80, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(80, 3).Value  = "Los Lagos"
$ws.Cells.Item(80, 4).Value  = 44511
$ws.Cells.Item(80, 5).Value  = 10
$ws.Cells.Item(80, 6).Value  = "Fruta"
$ws.Cells.Item(80, 7).Value  = 100108
$ws.Cells.Item(80, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(80, 9).Value  = 100108002
$ws.Cells.Item(80, 10).Value = "Mango"
$ws.Cells.Item(80, 11).Value = "Sin especificar"
$ws.Cells.Item(80, 12).Value = "Primera"
$ws.Cells.Item(80, 13).Value = 80
$ws.Cells.Item(80, 14).Value = 8000
$ws.Cells.Item(80, 15).Value = 8500
$ws.Cells.Item(80, 16).Value = 8250
$ws.Cells.Item(80, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(80, 18).Value = "Perú"
$ws.Cells.Item(80, 19).Value = 2062
$ws.Cells.Item(80, 20).Value = 4

# New record 2 -> row 81
$ws.Cells.Item(81, 1).Value  = 4
$ws.Cells.Item(81, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(81, 3).Value  = "Los Lagos"
$ws.Cells.Item(81, 4).Value  = 44511
$ws.Cells.Item(81, 5).Value  = 10
$ws.Cells.Item(81, 6).Value  = "Fruta"
$ws.Cells.Item(81, 7).Value  = 100108
$ws.Cells.Item(81, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(81, 9).Value  = 100108002
$ws.Cells.Item(81, 10).Value = "Mango"
$ws.Cells.Item(81, 11).Value = "Sin especificar"
$ws.Cells.Item(81, 12).Value = "Segunda"
$ws.Cells.Item(81, 13).Value = 20
$ws.Cells.Item(81, 14).Value = 6000
$ws.Cells.Item(81, 15).Value = 6000
$ws.Cells.Item(81, 16).Value = 6000
$ws.Cells.Item(81, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(81, 18).Value = "Perú"
$ws.Cells.Item(81, 19).Value = 1500
$ws.Cells.Item(81, 20).Value = 4
